## Updated cryptos list on Thu May 11 07:41:23 UTC 2023 with GitHub Actions
##
## Refreshes the Price (column D) and Volume(1h) (column E) figures on the
## crypto tracker sheet with the latest scraped values. Both columns hold
## plain text (not numbers/percent-formatted numerics) in this sheet, so
## every write below is forced to land as literal text -- matching how the
## sheet was originally populated -- instead of letting Excel's "looks like
## a number" auto-detection convert things like "312.45" into a numeric
## value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.570.18' },
    @{ Cell = 'E2'; Value = '  -0.14%  ' },
    @{ Cell = 'D3'; Value = '1.837.21' },
    @{ Cell = 'E3'; Value = '  -0.34%  ' },
    @{ Cell = 'E4'; Value = '  -0.15%  ' },
    @{ Cell = 'D5'; Value = '312.45' },
    @{ Cell = 'E6'; Value = '  -0.11%  ' },
    @{ Cell = 'D7'; Value = '0.4282' },
    @{ Cell = 'E7'; Value = '  -0.28%  ' },
    @{ Cell = 'D8'; Value = '0.3661' },
    @{ Cell = 'E8'; Value = '  +0.67%  ' },
    @{ Cell = 'D9'; Value = '0.07275' },
    @{ Cell = 'E9'; Value = '  -0.53%  ' },
    @{ Cell = 'D10'; Value = '0.8656' },
    @{ Cell = 'E10'; Value = '  -1.62%  ' },
    @{ Cell = 'D11'; Value = '20.71' },
    @{ Cell = 'E11'; Value = '  +0.42%  ' },
    @{ Cell = 'D12'; Value = '1.818.42' },
    @{ Cell = 'E12'; Value = '  -7.41%  ' },
    @{ Cell = 'D13'; Value = '5.471' },
    @{ Cell = 'E13'; Value = '  +2.40%  ' },
    @{ Cell = 'D14'; Value = '6.531' },
    @{ Cell = 'E14'; Value = '  +0.03%  ' },
    @{ Cell = 'D15'; Value = '0.06965' },
    @{ Cell = 'E15'; Value = '  +0.04%  ' },
    @{ Cell = 'E16'; Value = '  -0.12%  ' },
    @{ Cell = 'D17'; Value = '80.73' },
    @{ Cell = 'E17'; Value = '  +1.45%  ' },
    @{ Cell = 'D18'; Value = '0.000008894' },
    @{ Cell = 'E18'; Value = '  -0.91%  ' },
    @{ Cell = 'E19'; Value = '  -0.36%  ' },
    @{ Cell = 'E20'; Value = '  +0.46%  ' },
    @{ Cell = 'D21'; Value = '27.416.95' },
    @{ Cell = 'E21'; Value = '  -2.29%  ' },
    @{ Cell = 'D22'; Value = '5.163' },
    @{ Cell = 'E22'; Value = '  +3.57%  ' },
    @{ Cell = 'D23'; Value = '10.89' },
    @{ Cell = 'E23'; Value = '  +5.67%  ' },
    @{ Cell = 'D24'; Value = '2.031.69' },
    @{ Cell = 'E24'; Value = '  -4.25%  ' },
    @{ Cell = 'D25'; Value = '1.992' },
    @{ Cell = 'E25'; Value = '  +0.08%  ' },
    @{ Cell = 'D26'; Value = '154.85' },
    @{ Cell = 'E26'; Value = '  -0.38%  ' },
    @{ Cell = 'E27'; Value = '  +2.15%  ' },
    @{ Cell = 'D28'; Value = '5.171' },
    @{ Cell = 'E28'; Value = '  -0.81%  ' },
    @{ Cell = 'D29'; Value = '114.40' },
    @{ Cell = 'E29'; Value = '  -4.47%  ' },
    @{ Cell = 'D30'; Value = '1.824' },
    @{ Cell = 'E30'; Value = '  -3.14%  ' },
    @{ Cell = 'D31'; Value = '0.08862' },
    @{ Cell = 'E31'; Value = '  -0.45%  ' },
    @{ Cell = 'D32'; Value = '0.7512' },
    @{ Cell = 'E32'; Value = '  -1.02%  ' },
    @{ Cell = 'D33'; Value = '2.998' },
    @{ Cell = 'E33'; Value = '  +1.15%  ' },
    @{ Cell = 'D34'; Value = '4.547' },
    @{ Cell = 'E34'; Value = '  +0.55%  ' },
    @{ Cell = 'D35'; Value = '1.135' },
    @{ Cell = 'E35'; Value = '  +0.48%  ' },
    @{ Cell = 'E36'; Value = '  -0.08%  ' },
    @{ Cell = 'D37'; Value = '1.097' },
    @{ Cell = 'E37'; Value = '  -0.58%  ' },
    @{ Cell = 'D38'; Value = '0.05329' },
    @{ Cell = 'E38'; Value = '  -2.87%  ' },
    @{ Cell = 'D39'; Value = '0.01940' },
    @{ Cell = 'E39'; Value = '  +0.17%  ' },
    @{ Cell = 'D40'; Value = '2.798' },
    @{ Cell = 'E40'; Value = '  -1.04%  ' },
    @{ Cell = 'D41'; Value = '0.5081' },
    @{ Cell = 'E41'; Value = '  +0.11%  ' },
    @{ Cell = 'D42'; Value = '0.1651' },
    @{ Cell = 'E42'; Value = '  -0.88%  ' },
    @{ Cell = 'D43'; Value = '6.483' },
    @{ Cell = 'E43'; Value = '  -1.95%  ' },
    @{ Cell = 'E44'; Value = '  -0.77%  ' },
    @{ Cell = 'D45'; Value = '10.42' },
    @{ Cell = 'E45'; Value = '  +0.96%  ' },
    @{ Cell = 'D46'; Value = '105.58' },
    @{ Cell = 'E46'; Value = '  -0.58%  ' },
    @{ Cell = 'E47'; Value = '  -1.24%  ' },
    @{ Cell = 'D48'; Value = '0.4684' },
    @{ Cell = 'E48'; Value = '  +0.91%  ' },
    @{ Cell = 'E49'; Value = '  -0.15%  ' },
    @{ Cell = 'D50'; Value = '1.622' },
    @{ Cell = 'E50'; Value = '  -1.01%  ' },
    @{ Cell = 'D51'; Value = '1.742' },
    @{ Cell = 'E51'; Value = '  -0.48%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Value

    # Decide up front whether Excel's type-inference would swallow this
    # string as a real number (e.g. "312.45", "20.71", "0.000008894").
    # Values like "27.570.18" (two dots) or the "  -0.14%  " percent
    # strings never parse as a plain number and are safe to assign as-is.
    $looksNumeric = $text -match '^[+-]?[0-9]+(\.[0-9]+)?$'

    if ($looksNumeric) {
        # Force a literal-text write: set a text number format first so the
        # incoming string isn't coerced into a numeric cell value, then
        # restore the cell to the workbook's default "Normal" style so no
        # stray formatting is left behind (the source cells carry no
        # explicit style).
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}
